# Update "想去人数" (F column) figures and one cover image URL (I column)
# on the "展览" sheet and the "全部类型" sheet, reflecting a newer data pull.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F2").Value = 10119
$ws1.Range("F3").Value = 423
$ws1.Range("F5").Value = 21
$ws1.Range("F6").Value = 284
$ws1.Range("F8").Value = 479
$ws1.Range("F12").Value = 1047
$ws1.Range("F13").Value = 3146
$ws1.Range("F14").Value = 2353
$ws1.Range("F16").Value = 2075
$ws1.Range("F17").Value = 2075
$ws1.Range("F23").Value = 53
$ws1.Range("F25").Value = 6
$ws1.Range("F26").Value = 14
$ws1.Range("F32").Value = 577
$ws1.Range("F33").Value = 47
$ws1.Range("F34").Value = 229
$ws1.Range("F36").Value = 30
$ws1.Range("I36").Value = "//i1.hdslb.com/bfs/openplatform/202403/Kd0niodt1710905544733.jpeg"
$ws1.Range("F37").Value = 319
$ws1.Range("F38").Value = 1660
$ws1.Range("F39").Value = 110
$ws1.Range("F40").Value = 418
$ws1.Range("F42").Value = 436
$ws1.Range("F43").Value = 949
$ws1.Range("F45").Value = 348

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F2").Value = 10119
$ws4.Range("F3").Value = 423
$ws4.Range("F6").Value = 21
$ws4.Range("F8").Value = 284
$ws4.Range("F10").Value = 479
$ws4.Range("F13").Value = 1047
$ws4.Range("F14").Value = 3146
$ws4.Range("F15").Value = 2353
$ws4.Range("F16").Value = 2075
$ws4.Range("F17").Value = 2075
$ws4.Range("F23").Value = 53
$ws4.Range("F25").Value = 6
$ws4.Range("F26").Value = 14
$ws4.Range("F32").Value = 577
$ws4.Range("F36").Value = 47
$ws4.Range("F37").Value = 229
$ws4.Range("F39").Value = 30
$ws4.Range("I39").Value = "//i1.hdslb.com/bfs/openplatform/202403/Kd0niodt1710905544733.jpeg"
$ws4.Range("F41").Value = 319
$ws4.Range("F42").Value = 1660
$ws4.Range("F43").Value = 110
$ws4.Range("F45").Value = 418
$ws4.Range("F47").Value = 436
$ws4.Range("F48").Value = 949
$ws4.Range("F49").Value = 348
